# Automatic update of files.
# Bump the "Förändrad" (Changed) date in column C by one day for every
# data row (rows 2 through 39) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 39; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value2 = $cell.Value2 + 1
}
